$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-11, columns A-G
# A: Colaborador_id, B: Colaborador_nome, C: Departamento,
# D: Motivo_da_ausência, E: Horas_de_ausência, F: Data_da_ausência, G: Salário

$data = @(
    @{Row=2;  A=242;   B="Emanuel Fernandes";        C="TI";                   D="Outros";              E=5; F=45085; G=7443.82}
    @{Row=3;  A=50351; B="Benjamin Cardoso";          C="Jurídico";             D="Consulta médica";     E=3; F=45104; G=9879.280000000001}
    @{Row=4;  A=90122; B="Srta. Ana Júlia Castro";    C="Operações";            D="Viagem de negócios";  E=2; F=45090; G=9746.450000000001}
    @{Row=5;  A=48076; B="Sra. Ana Clara Farias";     C="Operações";            D="Consulta médica";     E=1; F=45083; G=6315.36}
    @{Row=6;  A=62562; B="Breno Rocha";               C="Vendas";               D="Outros";              E=7; F=45097; G=4760.37}
    @{Row=7;  A=84758; B="Luiz Miguel Porto";         C="Operações";            D="Problemas pessoais";  E=3; F=45097; G=9450.4}
    @{Row=8;  A=57868; B="Dra. Gabrielly Gomes";      C="Marketing";            D="Problemas pessoais";  E=3; F=45079; G=10420.33}
    @{Row=9;  A=96191; B="Miguel da Paz";             C="TI";                   D="Outros";              E=3; F=45106; G=7918.4}
    @{Row=10; A=54496; B="Srta. Rafaela Carvalho";    C="Operações";            D="Doença";              E=4; F=45083; G=8775.34}
    @{Row=11; A=78845; B="Vitor Gabriel Pereira";     C="Financeiro";           D="Viagem de negócios";  E=5; F=45083; G=6127.13}
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value = $rowData.A
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
}
